$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24; this shifts rows 24..111 down to 25..112
$ws.Rows("24:24").Insert()

# Populate the newly inserted row 24 with its data
$ws.Cells.Item(24, 1).Value = 3
$ws.Cells.Item(24, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(24, 3).Value = "Coquimbo"
$ws.Cells.Item(24, 4).Value = 44565
$ws.Cells.Item(24, 5).Value = 5
$ws.Cells.Item(24, 6).Value = 100112052
$ws.Cells.Item(24, 7).Value = "Albahaca"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 140
$ws.Cells.Item(24, 11).Value = 5000
$ws.Cells.Item(24, 12).Value = 5500
$ws.Cells.Item(24, 13).Value = 5286
$ws.Cells.Item(24, 14).Value = "`$/docena de matas"
$ws.Cells.Item(24, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(24, 16).Value = 881
$ws.Cells.Item(24, 17).Value = 6
$ws.Cells.Item(24, 18).Value = "Hortaliza"
